# Updates crypto price (D) and 1h volume-change (E) columns to the latest
# scraped snapshot. Price cells that look numeric must stay text (the sheet
# stores "Price" as plain strings, e.g. "1.00" / "10.70" / "0.0170"), so those
# are forced to Text via NumberFormat "@" before the write and restored to the
# default "Normal" style afterwards so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "28.010.96"
    "E2" = "  +1.12%  "
    "D3" = "1.641.84"
    "E3" = "  +0.53%  "
    "E4" = "  -0.05%  "
    "D5" = "212.97"
    "E5" = "  +0.36%  "
    "E6" = "  +0.45%  "
    "D7" = "1.00"
    "E7" = "  -0.09%  "
    "D8" = "23.63"
    "E8" = "  +1.89%  "
    "E9" = "  -1.89%  "
    "E10" = "  +0.42%  "
    "D11" = "0.0883"
    "E11" = "  +2.57%  "
    "D12" = "1.874.03"
    "E12" = "  +0.48%  "
    "D13" = "1.641.94"
    "E13" = "  +0.58%  "
    "D14" = "0.575"
    "E14" = "  +3.65%  "
    "E15" = "  +1.40%  "
    "D16" = "65.94"
    "E16" = "  +1.16%  "
    "D17" = "28.004.97"
    "E17" = "  +1.21%  "
    "D18" = "236.54"
    "E18" = "  +2.91%  "
    "E19" = "  +0.55%  "
    "D20" = "7.64"
    "E20" = "  +0.99%  "
    "D21" = "1.00"
    "E21" = "  -0.05%  "
    "D22" = "10.70"
    "E22" = "  +0.37%  "
    "E23" = "  +0.89%  "
    "D24" = "2.10"
    "E24" = "  -1.81%  "
    "D25" = "151.57"
    "E25" = "  +1.85%  "
    "D26" = "6.97"
    "E26" = "  +1.45%  "
    "D27" = "15.71"
    "E27" = "  +0.88%  "
    "E28" = "  +0.18%  "
    "E30" = "  +0.37%  "
    "E31" = "  +0.56%  "
    "E32" = "  +2.04%  "
    "D33" = "3.13"
    "E33" = "  +1.58%  "
    "D34" = "1.419.11"
    "E35" = "  +2.68%  "
    "D36" = "2.36"
    "E36" = "  +1.33%  "
    "D37" = "0.0170"
    "E37" = "  +1.69%  "
    "D38" = "0.884"
    "E38" = "  +0.87%  "
    "E39" = "  -0.02%  "
    "D40" = "0.901"
    "E40" = "  -4.84%  "
    "E41" = "  +0.99%  "
    "E43" = "  +6.78%  "
    "D44" = "66.71"
    "E44" = "  -1.68%  "
    "D45" = "5.51"
    "E45" = "  +2.85%  "
    "E46" = "  -0.12%  "
    "D47" = "1.783.14"
    "E47" = "  +0.62%  "
    "D48" = "87.89"
    "E48" = "  +0.36%  "
    "E49" = "  +0.85%  "
    "E50" = "  +0.39%  "
    "D51" = "7.64"
    "E51" = "  -0.80%  "
}

$forceText = @(
    "D5"
    "D7"
    "D8"
    "D11"
    "D14"
    "D16"
    "D18"
    "D20"
    "D21"
    "D22"
    "D24"
    "D25"
    "D26"
    "D27"
    "D33"
    "D36"
    "D37"
    "D38"
    "D40"
    "D44"
    "D45"
    "D48"
    "D51"
)

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    if ($forceText -contains $cell) {
        $range.NumberFormat = "@"
        $range.Value = $updates[$cell]
        $range.Style = "Normal"
    } else {
        $range.Value = $updates[$cell]
    }
}
